# Commit: "added power analysis and started website rework"
#
# For this particular workbook (docs/data/survey-2.xlsx) the substantive
# change is a corrected / rephrased source reference for the "Climate
# Crisis Behavior" items: the old citation "(Sorrell et al., 2020)" used
# in the Reference column (E3:E5) is replaced with "Based on Sorrel et
# al., 2020". (Row 6 / E6 has no reference cell, same as before.)
#
# All of the other bytes shown in the raw xml diff (fileVersion,
# rupBuild, calcId, xr:* revision GUIDs, the x15ac:absPath user folder,
# window/selection coordinates, row heights / dyDescent, styles.xml
# namespace list, ...) are incidental Excel-version/environment save
# artifacts, not authored content changes, so they are intentionally not
# reproduced here.

$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Tabelle1") {
        $ws = $sheet
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$newReference = "Based on Sorrel et al., 2020"

$ws.Range("E3").Value = $newReference
$ws.Range("E4").Value = $newReference
$ws.Range("E5").Value = $newReference
